$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.206.46"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.110.66"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.41%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.83"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.19%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.520"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.89%  "
$ws.Range("E9").Value = "  +0.11%  "
$ws.Range("E10").Value = "  -1.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.477"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.69%  "
$ws.Range("E12").Value = "  -1.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "36.66"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.80%  "
$ws.Range("E14").Value = "  -1.96%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.627.40"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.153.57"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.07%  "
$ws.Range("E17").Value = "  -1.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.109.38"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.57"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "490.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.83"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.699"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "83.83"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.10%  "
$ws.Range("E25").Value = "  -2.56%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.54"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.28%  "
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.89"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.65%  "
$ws.Range("E29").Value = "  -2.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.67"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.03%  "
$ws.Range("E31").Value = "  -2.37%  "
$ws.Range("E32").Value = "  -1.40%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0₃0938"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.21%  "
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.77"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.86%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.969"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "46.77"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.39%  "
$ws.Range("E38").Value = "  -4.34%  "
$ws.Range("E39").Value = "  +0.65%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.307"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.46"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.46%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "386.11"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.798.09"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.98%  "
$ws.Range("E44").Value = "  -8.67%  "
$ws.Range("E45").Value = "  -2.62%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "135.10"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.19%  "
$ws.Range("E47").Value = "  +0.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "24.99"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.59%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.19"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.83%  "
$ws.Range("E50").Value = "  -1.85%  "
$ws.Range("E51").Value = "  -1.98%  "
